$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-23 changes from serial date 45221 (2023-10-22)
# to serial date 45224 (2023-10-25).
for ($row = 2; $row -le 23; $row++) {
    $ws.Cells.Item($row, 3).Value = 45224
}
